$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.2579443333333333"
$ws.Range("H2").Value = [double]"0.773833"
$ws.Range("I2").Value = [double]"0.05467096027587241"
$ws.Range("J2").Value = [double]"0.0546709602758724"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.3541136666666667"
$ws.Range("N2").Value = [double]"1.062341"
$ws.Range("O2").Value = [double]"0.002054893867096745"
$ws.Range("P2").Value = [double]"0.002054893867096745"
$ws.Range("Q2").Value = [double]"0.09134161367255555"
$ws.Range("R2").Value = [double]"0.822074523053"
$ws.Range("S2").Value = [double]"0.00011234302097918"
$ws.Range("T2").Value = [double]"0.00011234302097918"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.2579443333333333"
$ws.Range("H3").Value = [double]"0.773833"
$ws.Range("I3").Value = [double]"0.05467096027587241"
$ws.Range("J3").Value = [double]"0.0546709602758724"
$ws.Range("K3").Value = [double]"1"
$ws.Range("L3").Value = [double]"0.3333333333333333"
$ws.Range("M3").Value = [double]"0.03021466666666667"
$ws.Range("N3").Value = [double]"0.090644"
$ws.Range("O3").Value = [double]"0.0001753333437089572"
$ws.Range("P3").Value = [double]"0.0001753333437089573"
$ws.Range("Q3").Value = [double]"0.007793702050222223"
$ws.Range("R3").Value = [double]"0.070143318452"
$ws.Range("S3").Value = [double]"9.585642268948285E-06"
$ws.Range("T3").Value = [double]"9.585642268948285E-06"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.2579443333333333"
$ws.Range("H4").Value = [double]"0.773833"
$ws.Range("I4").Value = [double]"0.05467096027587241"
$ws.Range("J4").Value = [double]"0.0546709602758724"
$ws.Range("M4").Value = [double]"171.9426576666667"
$ws.Range("N4").Value = [double]"515.827973"
$ws.Range("O4").Value = [double]"0.9977697727891942"
$ws.Range("P4").Value = [double]"0.9977697727891944"
$ws.Range("Q4").Value = [double]"44.35163420338989"
$ws.Range("R4").Value = [double]"399.164707830509"
$ws.Range("S4").Value = [double]"0.05454903161262428"
$ws.Range("T4").Value = [double]"0.05454903161262428"
$ws.Range("I5").Value = [double]"0.4403936734732808"
$ws.Range("J5").Value = [double]"0.4403936734732807"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.3541136666666667"
$ws.Range("N5").Value = [double]"1.062341"
$ws.Range("O5").Value = [double]"0.002054893867096745"
$ws.Range("P5").Value = [double]"0.002054893867096745"
$ws.Range("Q5").Value = [double]"0.7357885901994445"
$ws.Range("R5").Value = [double]"6.622097311795001"
$ws.Range("S5").Value = [double]"0.000904962258728451"
$ws.Range("T5").Value = [double]"0.000904962258728451"
$ws.Range("I6").Value = [double]"0.4403936734732808"
$ws.Range("J6").Value = [double]"0.4403936734732807"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.3333333333333333"
$ws.Range("M6").Value = [double]"0.03021466666666667"
$ws.Range("N6").Value = [double]"0.090644"
$ws.Range("O6").Value = [double]"0.0001753333437089572"
$ws.Range("P6").Value = [double]"0.0001753333437089573"
$ws.Range("Q6").Value = [double]"0.06278099119777779"
$ws.Range("R6").Value = [double]"0.5650289207800001"
$ws.Range("S6").Value = [double]"7.721569531834103E-05"
$ws.Range("T6").Value = [double]"7.721569531834103E-05"
$ws.Range("I7").Value = [double]"0.4403936734732808"
$ws.Range("J7").Value = [double]"0.4403936734732807"
$ws.Range("M7").Value = [double]"171.9426576666667"
$ws.Range("N7").Value = [double]"515.827973"
$ws.Range("O7").Value = [double]"0.9977697727891942"
$ws.Range("P7").Value = [double]"0.9977697727891944"
$ws.Range("Q7").Value = [double]"357.2678989506262"
$ws.Range("R7").Value = [double]"3215.411090555635"
$ws.Range("S7").Value = [double]"0.439411495519234"
$ws.Range("T7").Value = [double]"0.439411495519234"
$ws.Range("G8").Value = [double]"2.217259"
$ws.Range("H8").Value = [double]"6.651777"
$ws.Range("I8").Value = [double]"0.4699451123575263"
$ws.Range("J8").Value = [double]"0.4699451123575263"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.3541136666666667"
$ws.Range("N8").Value = [double]"1.062341"
$ws.Range("O8").Value = [double]"0.002054893867096745"
$ws.Range("P8").Value = [double]"0.002054893867096745"
$ws.Range("Q8").Value = [double]"0.7851617144396666"
$ws.Range("R8").Value = [double]"7.066455429957"
$ws.Range("S8").Value = [double]"0.0009656873292555714"
$ws.Range("T8").Value = [double]"0.0009656873292555716"
$ws.Range("G9").Value = [double]"2.217259"
$ws.Range("H9").Value = [double]"6.651777"
$ws.Range("I9").Value = [double]"0.4699451123575263"
$ws.Range("J9").Value = [double]"0.4699451123575263"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.03021466666666667"
$ws.Range("N9").Value = [double]"0.090644"
$ws.Range("O9").Value = [double]"0.0001753333437089572"
$ws.Range("P9").Value = [double]"0.0001753333437089573"
$ws.Range("Q9").Value = [double]"0.06699374159866667"
$ws.Range("R9").Value = [double]"0.602943674388"
$ws.Range("S9").Value = [double]"8.23970479093267E-05"
$ws.Range("T9").Value = [double]"8.239704790932671E-05"
$ws.Range("G10").Value = [double]"2.217259"
$ws.Range("H10").Value = [double]"6.651777"
$ws.Range("I10").Value = [double]"0.4699451123575263"
$ws.Range("J10").Value = [double]"0.4699451123575263"
$ws.Range("M10").Value = [double]"171.9426576666667"
$ws.Range("N10").Value = [double]"515.827973"
$ws.Range("O10").Value = [double]"0.9977697727891942"
$ws.Range("P10").Value = [double]"0.9977697727891944"
$ws.Range("Q10").Value = [double]"381.2414051953357"
$ws.Range("R10").Value = [double]"3431.172646758022"
$ws.Range("S10").Value = [double]"0.4688970279803614"
$ws.Range("T10").Value = [double]"0.4688970279803615"
$ws.Range("E11").Value = [double]"2"
$ws.Range("F11").Value = [double]"0.6666666666666666"
$ws.Range("G11").Value = [double]"0.1650883333333333"
$ws.Range("H11").Value = [double]"0.495265"
$ws.Range("I11").Value = [double]"0.03499025389332058"
$ws.Range("J11").Value = [double]"0.03499025389332058"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.3541136666666667"
$ws.Range("N11").Value = [double]"1.062341"
$ws.Range("O11").Value = [double]"0.002054893867096745"
$ws.Range("P11").Value = [double]"0.002054893867096745"
$ws.Range("Q11").Value = [double]"0.05846003504055555"
$ws.Range("R11").Value = [double]"0.5261403153650001"
$ws.Range("S11").Value = [double]"7.190125813354246E-05"
$ws.Range("T11").Value = [double]"7.190125813354246E-05"
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"0.1650883333333333"
$ws.Range("H12").Value = [double]"0.495265"
$ws.Range("I12").Value = [double]"0.03499025389332058"
$ws.Range("J12").Value = [double]"0.03499025389332058"
$ws.Range("K12").Value = [double]"1"
$ws.Range("L12").Value = [double]"0.3333333333333333"
$ws.Range("M12").Value = [double]"0.03021466666666667"
$ws.Range("N12").Value = [double]"0.090644"
$ws.Range("O12").Value = [double]"0.0001753333437089572"
$ws.Range("P12").Value = [double]"0.0001753333437089573"
$ws.Range("Q12").Value = [double]"0.004988088962222223"
$ws.Range("R12").Value = [double]"0.04489280066"
$ws.Range("S12").Value = [double]"6.134958212341257E-06"
$ws.Range("T12").Value = [double]"6.134958212341257E-06"
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"0.1650883333333333"
$ws.Range("H13").Value = [double]"0.495265"
$ws.Range("I13").Value = [double]"0.03499025389332058"
$ws.Range("J13").Value = [double]"0.03499025389332058"
$ws.Range("M13").Value = [double]"171.9426576666667"
$ws.Range("N13").Value = [double]"515.827973"
$ws.Range("O13").Value = [double]"0.9977697727891942"
$ws.Range("P13").Value = [double]"0.9977697727891944"
$ws.Range("Q13").Value = [double]"28.38572678309389"
$ws.Range("R13").Value = [double]"255.471541047845"
$ws.Range("S13").Value = [double]"0.03491221767697469"
$ws.Range("T13").Value = [double]"0.03491221767697469"
